$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we are about to write keep their literal text representation
# (Excel auto-converts numeric-looking strings to numbers otherwise, which would
# drop things like trailing zeros or change multi-dot "thousands" style numbers).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.025.28"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.862.44"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "312.27"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "0.5096"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("D8").Value = "0.3847"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "0.08295"
$ws.Range("E9").Value = "  -7.83%  "
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "41.53"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "6.235"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "20.60"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").Value = "1.862.99"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "7.236"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "90.88"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "0.06632"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "6.043"
$ws.Range("D23").Value = "28.052.18"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").Value = "2.233"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "3.397"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.538"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.077.15"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "157.89"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "20.53"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "124.93"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.1056"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "1.037"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.886"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "3.595"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "9.412"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "0.06541"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02423"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2174"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.204"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6469"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").Value = "4.990"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.226"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "11.19"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6128"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.16"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "1.287"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "3.656"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.016"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.210"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "120.17"
$ws.Range("E51").Value = "  -0.56%  "
